$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new worksheet "getCacheKeyAndValue" right after "getAllCacheNames"
# (i.e. at the end of the tab strip), mirroring that sheet's layout but with
# two extra columns (cacheName / entityName) inserted between "description"
# and "rspStatus".
# ---------------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("getAllCacheNames")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "getCacheKeyAndValue"

# Copy the header row formatting (bold/filled/bordered header style) from the
# existing "getAllCacheNames" sheet onto the new header row.
$srcSheet.Range("A1:F1").Copy()
$newSheet.Range("A1:G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column widths matching the target layout.
$newSheet.Columns.Item(1).ColumnWidth = 27.857142857142854
$newSheet.Columns.Item(2).ColumnWidth = 44.714285714285715
$newSheet.Columns.Item(3).ColumnWidth = 26.57142857142857
$newSheet.Columns.Item(4).ColumnWidth = 13.428571428571429

# Write the new-string-introducing cells first, in the exact order the
# strings first appear, so they line up with how the workbook records them.
$newSheet.Range("A2").Value = "jinzu-connector-cache-test2"
$newSheet.Range("D1").Value = "entityName"
$newSheet.Range("D2").Value = "Site"
$newSheet.Range("C1").Value = "cacheName"
$newSheet.Range("C2").Value = "data"
$newSheet.Range("A3").Value = "jinzu-connector-cache-test3"
$newSheet.Range("D3").Value = "Inverter"
$newSheet.Range("C3").Value = "datasource_auth"
$newSheet.Range("B2").Value = " data,check get cache key and value"
$newSheet.Range("B3").Value = "datasource_auth,check get cache key and value"

# Remaining header cells reuse already-existing shared strings.
$newSheet.Range("A1").Value = "test-id"
$newSheet.Range("B1").Value = "description"
$newSheet.Range("E1").Value = "rspStatus"
$newSheet.Range("F1").Value = "rspCode"
$newSheet.Range("G1").Value = "rspMessage"

# Page setup / margins to match the other data sheets.
$newSheet.PageSetup.Orientation = 1
$newSheet.PageSetup.LeftMargin = $excel.InchesToPoints(0.7)
$newSheet.PageSetup.RightMargin = $excel.InchesToPoints(0.7)
$newSheet.PageSetup.TopMargin = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.BottomMargin = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.HeaderMargin = $excel.InchesToPoints(0.3)
$newSheet.PageSetup.FooterMargin = $excel.InchesToPoints(0.3)

# Selection on the new sheet.
$newSheet.Range("B10").Select()

# ---------------------------------------------------------------------------
# Update the view state of "getAllCacheNames": it is no longer the selected
# tab, and its remembered selection moves to G5.
# ---------------------------------------------------------------------------
$srcSheet.Range("G5").Select()

# Activating the new sheet last makes it the active tab (tabSelected="1" /
# activeTab) while leaving the other sheet's selection as set above.
$newSheet.Activate()
